# Weekly update: insert a new observation row for the latest week at the
# top of the data block (row 279), pushing the existing historical rows
# down by one. This mirrors the source system prepending the newest
# weekly price record to the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 279:351 down to 280:352, creating a blank row 279.
$ws.Rows.Item(279).Insert()

# Populate the newly inserted row 279 with the new weekly record.
$ws.Cells.Item(279, 1).Value  = 6
$ws.Cells.Item(279, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(279, 3).Value  = "Metropolitana"
$ws.Cells.Item(279, 4).Value  = 45135
$ws.Cells.Item(279, 5).Value  = 13
$ws.Cells.Item(279, 6).Value  = 100112029
$ws.Cells.Item(279, 7).Value  = "Orégano"
$ws.Cells.Item(279, 8).Value  = "Sin especificar"
$ws.Cells.Item(279, 9).Value  = "Primera"
$ws.Cells.Item(279, 10).Value = 30
$ws.Cells.Item(279, 11).Value = 20000
$ws.Cells.Item(279, 12).Value = 20000
$ws.Cells.Item(279, 13).Value = 20000
$ws.Cells.Item(279, 14).Value = "$/docena de atados"
$ws.Cells.Item(279, 15).Value = "Región Metropolitana"
$ws.Cells.Item(279, 16).Value = 6667
$ws.Cells.Item(279, 17).Value = 3
$ws.Cells.Item(279, 18).Value = "Hortaliza"
